$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains its text formatting so values
# like "1.004" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.044.98'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '1.928.48'
$ws.Range('E3').Value = '  +1.35%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '325.32'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = '0.4592'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.3823'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '0.07763'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').Value = '0.9806'
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').Value = '22.62'
$ws.Range('E11').Value = '  +2.63%  '
$ws.Range('D12').Value = '1.919.61'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').Value = '5.718'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').Value = '6.984'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '0.06999'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '84.89'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').Value = '1.005'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '0.000009502'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '16.75'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = '29.048.62'
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').Value = '11.10'
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('D24').Value = '2.155.94'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('D25').Value = '2.055'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').Value = '158.28'
$ws.Range('E26').Value = '  +0.94%  '
$ws.Range('D27').Value = '19.11'
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('D28').Value = '5.637'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').Value = '117.70'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '1.845'
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('D31').Value = '0.09325'
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('D32').Value = '0.8653'
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('D33').Value = '5.114'
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('D34').Value = '1.249'
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('D35').Value = '3.014'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '0.05709'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D38').Value = '1.003'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').Value = '0.02057'
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('D40').Value = '3.096'
$ws.Range('D41').Value = '7.472'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '0.5520'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').Value = '0.1761'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').Value = '9.381'
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('D45').Value = '0.000002808'
$ws.Range('E45').Value = '  +10.69%  '
$ws.Range('D46').Value = '2.185'
$ws.Range('E46').Value = '  +4.61%  '
$ws.Range('D47').Value = '0.5183'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('D48').Value = '0.06935'
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').Value = '11.20'
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('D50').Value = '111.07'
$ws.Range('E50').Value = '  -0.29%  '
$ws.Range('D51').Value = '1.770'
$ws.Range('E51').Value = '  -0.19%  '
